# Ajuste manejo de alertas
# - Update the stored password hint value in A2.
# - Remove the (alert-style) underline formatting that was applied to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Willian00020"
$ws.Range("B2").Font.Underline = -4142  # xlUnderlineStyleNone
